$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "57.959.07"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.85%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.073.29"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.82%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "514.70"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.75%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "140.88"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +1.54%  "

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.02%  "

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.79%  "

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +2.34%  "

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.13%  "

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +3.24%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "3.599.96"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.64%  "

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +1.57%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.65"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +6.31%  "

$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.68%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "57.991.15"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +1.75%  "

$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value = "Polkadot"
$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.16"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +4.52%  "

$ws.Range("B18").NumberFormat = "@"
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").NumberFormat = "@"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.068.17"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.48%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.83"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -1.49%  "

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.14%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "333.10"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.16%  "

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.00%  "

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.08%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.09"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.02%  "

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +3.52%  "

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.37%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0₃0904"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -2.66%  "

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.41%  "

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +4.24%  "

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.30%  "

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +3.28%  "

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.74%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "155.24"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.68%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.54"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +1.37%  "

$ws.Range("B35").NumberFormat = "@"
$ws.Range("B35").Value = "EnergySwap"
$ws.Range("C35").NumberFormat = "@"
$ws.Range("C35").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "27.31"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +4.03%  "

$ws.Range("B36").NumberFormat = "@"
$ws.Range("B36").Value = "Aptos"
$ws.Range("C36").NumberFormat = "@"
$ws.Range("C36").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.99"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +3.12%  "

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +3.86%  "

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +1.36%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.114.10"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.85%  "

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +2.22%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "36.58"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.61%  "

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +0.08%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.655"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -1.88%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.278.55"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +2.20%  "

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +1.64%  "

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +2.02%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "20.50"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +4.40%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.942"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +1.40%  "

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +2.21%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.726"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +6.74%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "256.65"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +10.33%  "
